$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.388.73'
$ws.Range('E2').Value = '  -4.70%  '
$ws.Range('D3').Value = '3.355.95'
$ws.Range('E3').Value = '  -1.91%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '567.60'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.19'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.83%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '3.354.79'
$ws.Range('E8').Value = '  -1.94%  '
$ws.Range('E9').Value = '  -1.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.49'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.17%  '
$ws.Range('E11').Value = '  -2.58%  '
$ws.Range('E12').Value = '  -0.49%  '
$ws.Range('D13').Value = '3.924.76'
$ws.Range('E13').Value = '  -1.94%  '
$ws.Range('E14').Value = '  -0.14%  '
$ws.Range('D15').Value = '3.354.12'
$ws.Range('E16').Value = '  -2.23%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '24.88'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.68%  '
$ws.Range('D18').Value = '60.463.51'
$ws.Range('E18').Value = '  -4.60%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.62'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.32%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '9.35'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.48%  '
$ws.Range('B21').Value = 'Polkadot'
$ws.Range('C21').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.73'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.63%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '368.07'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.87%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.562'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('D24').Value = '3.484.42'
$ws.Range('E24').Value = '  -2.12%  '
$ws.Range('E25').Value = '  -0.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '69.81'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000113'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.63%  '
$ws.Range('E28').Value = '  +17.91%  '
$ws.Range('E29').Value = '  +7.88%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.994'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.06'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.155'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.69%  '
$ws.Range('E33').Value = '  -1.90%  '
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('D35').Value = '3.386.50'
$ws.Range('E35').Value = '  -1.88%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '23.03'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.88%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.37'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.95'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.28%  '
$ws.Range('E39').Value = '  +1.56%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '158.94'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.87%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0777'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.28%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.998'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.41'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.66%  '
$ws.Range('E44').Value = '  +10.44%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '41.01'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.41%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.752'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.99%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '23.98'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.69%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.60'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.42%  '
$ws.Range('E49').Value = '  +2.08%  '
$ws.Range('E50').Value = '  +12.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.897'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.98%  '
